$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("A1").Value = "forces"
$ws.Range("B1").Value = "canmove?"

# Column A data / formulas
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 0
$ws.Range("A4").Value = 0
$ws.Range("A5").Value = 0
$ws.Range("A6").Formula = "=100000*COS(2*PI()*60/360)"
$ws.Range("A7").Formula = "=-100000*SIN(2*PI()*60/360)"
$ws.Range("A8").Value = 0
$ws.Range("A9").Value = 0

# Column B data
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = 1
$ws.Range("B7").Value = 1
$ws.Range("B8").Value = 0
$ws.Range("B9").Value = 0

# Borders: thick rule under the header row
$ws.Range("A1:B1").Borders.Item(9).Weight = -4138

# Borders: box down the left side of column B (data column), with a
# medium cap along the top of the box under the header
$ws.Range("B2").Borders.Item(7).Weight = 2
$ws.Range("B2").Borders.Item(8).Weight = -4138
$ws.Range("B3:B9").Borders.Item(7).Weight = 2

# Selection / view
[void]$ws.Range("A16").Select()

# Page layout
$ws.PageSetup.Orientation = 1
